# "Generate Report for Handoff"
#
# The report moves from "In Translation" to "Ready for handoff" and the
# handoff timestamps are refreshed. Because the sheets share their string
# table, updating the cells below updates every occurrence of the old text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# Refresh the handoff timestamps
$wsOverview.Range("G2").Value = "2016-08-29 09:02:58"
$wsDeDe.Range("H2").Value     = "2016-08-29 09:02:58"
$wsZhCn.Range("H2").Value     = "2016-08-29 09:02:54"

# The longer "Ready for handoff" label no longer fits the old column width,
# so the Status columns get widened (mirrors Excel's own autofit reflow).
$wsOverview.Columns("E:F").ColumnWidth = 16.333333333333336
$wsZhCn.Columns("C:C").ColumnWidth     = 16.333333333333336
$wsDeDe.Columns("C:C").ColumnWidth     = 16.333333333333336
